# Auto-generated: applies quarterly "cryptos list" refresh (GitHub Actions bot).
# D column = Price (kept as text via a leading apostrophe so Excel does not
# reinterpret values like "0.0630" as the number 0.063), E column = Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.901.83"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").Value = "'1.729.84"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +3.45%  '
$ws.Range("E4").Value = '  -0.38%  '
$ws.Range("D5").Value = "'217.51"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("E6").Value = '  +1.36%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = "'23.81"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +11.21%  '
$ws.Range("E9").Value = '  +4.41%  '
$ws.Range("D10").Value = "'0.0630"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  +1.45%  '
$ws.Range("D11").Value = "'0.0899"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("D12").Value = "'1.972.53"
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("D13").Value = "'1.730.81"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  +4.03%  '
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("E15").Value = '  +6.41%  '
$ws.Range("E16").Value = '  +3.06%  '
$ws.Range("D17").Value = "'27.918.52"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("D18").Value = "'241.98"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +3.19%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").Value = "'0.0₃0750"
$ws.Range("D20").Style = 'Normal'
$ws.Range("E20").Value = '  +2.14%  '
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("E22").Value = '  +4.02%  '
$ws.Range("E23").Value = '  +5.11%  '
$ws.Range("E24").Value = '  +0.68%  '
$ws.Range("D25").Value = "'148.74"
$ws.Range("D25").Style = 'Normal'
$ws.Range("E25").Value = '  +1.49%  '
$ws.Range("E26").Value = '  +4.18%  '
$ws.Range("D27").Value = "'16.65"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +1.70%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("E30").Value = '  +1.61%  '
$ws.Range("E31").Value = '  +1.82%  '
$ws.Range("E32").Value = '  +2.92%  '
$ws.Range("E33").Value = '  +4.39%  '
$ws.Range("D34").Value = "'1.487.13"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -3.04%  '
$ws.Range("D35").Value = "'1.69"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -0.78%  '
$ws.Range("E36").Value = '  +6.55%  '
$ws.Range("E37").Value = '  +3.85%  '
$ws.Range("D38").Value = "'2.40"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  +2.36%  '
$ws.Range("D41").Value = "'72.12"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +6.74%  '
$ws.Range("D42").Value = "'5.87"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +6.17%  '
$ws.Range("E43").Value = '  -0.32%  '
$ws.Range("D44").Value = "'1.878.20"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  +3.48%  '
$ws.Range("D45").Value = "'2.29"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("D46").Value = "'0.792"
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("E47").Value = '  +10.08%  '
$ws.Range("D48").Value = "'91.86"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +1.68%  '
$ws.Range("E49").Value = '  +4.71%  '

# Rows 50/51: coin list order changed -- EnergySwap now ranks above Algorand,
# each carrying its own refreshed price/volume figures.
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'8.33"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +4.15%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = "'0.106"
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +2.52%  '
